# Automatische test-sync: 2025-07-22 12:35:50
# Appends the new "Testmail #11" row to the Logs sheet and updates the
# Dashboard summary count for the "Retour / Terugbetaling" category.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append row 11 -------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(11, 1).Value = "Mijn product is beschadigd geleverd."
$logs.Cells.Item(11, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(11, 3).Value = "Testmail #11: Mijn product is beschadigd geleverd."
$logs.Cells.Item(11, 4).Value = "Retour / Terugbetaling"
$logs.Cells.Item(11, 5).Value = "Beste klant,`nBedankt voor het melden van dit probleem. Om u beter van dienst te kunnen zijn, hebben wij wat meer informatie nodig over de beschadiging. Kunt u alstublieft een foto van het beschadigde product meesturen? Dit helpt ons om het probleem beter te begrijpen en een passende oplossing voor u te vinden.`nMet vriendelijke groet,`n[E-mailassistent] van [Bedrijfsnaam]"
$logs.Cells.Item(11, 6).Value = "2025-07-22 12:35:06"
$logs.Cells.Item(11, 7).Value = "Ja"
$logs.Cells.Item(11, 8).Value = "Nee"
$logs.Cells.Item(11, 9).Value = "Ja"
$logs.Cells.Item(11, 10).Value = "Ja"

# Keep the default row height (matches the other data rows, which also
# hold multi-line replies without a custom row height).
$logs.Rows.Item(11).AutoFit()

# --- Logs sheet: extend conditional formatting ranges to include row 11 --
$logs.Range("D2:D10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D11"))
$logs.Range("G2:G10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G11"))
$logs.Range("H2:H10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H11"))
$logs.Range("I2:I10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I11"))
$logs.Range("J2:J10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J11"))

# --- Dashboard sheet: bump the "Retour / Terugbetaling" count ------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(3, 2).Value = 4
